$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: rename label and keep rest the same
$ws.Range("A2").Value = "MSFT - Income Statement URL"

# Row 3: Balance Sheet
$ws.Range("A3").Value = "MSFT - Balance Sheet URL"
$ws.Range("C3").Value = $true
$ws.Range("D3").Value = "MarketData"
$ws.Range("E3").Value = "#sfcontent > div.rf_ctlwrap > div.rf_ctl2_opt > div.exportButton > span > a"
$ws.Range("G3").Value = "http://financials.morningstar.com/balance-sheet/bs.html?t=MSFT&region=usa&culture=en-US"

# Row 4: Ratios
$ws.Range("A4").Value = "MSFT - Ratios URL"
$ws.Range("C4").Value = $true
$ws.Range("D4").Value = "MarketData"
$ws.Range("E4").Value = "#financials > div.r_tbar0.positionrelative > div > a > div"
$ws.Range("G4").Value = "http://financials.morningstar.com/ratios/r.html?t=MSFT&region=usa&culture=en-US"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 22.7109375
$ws.Columns.Item(5).ColumnWidth = 75.140625

# Selection
[void]$ws.Range("B2").Select()
